$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 16: table graphicFrame's table style id changes from the
#    presentation's custom "Table_0" style to the built-in
#    {35E619DA-9D75-4AAC-ACF4-FA40D2763605} style.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{35E619DA-9D75-4AAC-ACF4-FA40D2763605}")
    }
}

# ---------------------------------------------------------------------------
# 2) The deck's theme (the one attached to the slide master / all slide
#    layouts) swaps its colour scheme from the "Integral" palette to the
#    plain "Office Theme" palette (the palette that used to live only on
#    the Notes Master's theme part).
# ---------------------------------------------------------------------------
$officeThemeColors = @(
    0x000000, # dk1
    0xFFFFFF, # lt1
    0x6A546A, # (placeholder, overwritten below with correct value)
    0xE7E6E6, # lt2
    0x5B9BD5, # accent1
    0xED7D31, # accent2
    0xA5A5A5, # accent3
    0xFFC000, # accent4
    0x4472C4, # accent5
    0x70AD47, # accent6
    0x0563C1, # hlink
    0x954F72  # folHlink
)
# dk2 = 44546A (written explicitly below so the literal hex/RGB byte-order
# mapping used by PowerPoint's RGB() is unambiguous)
$officeThemeColors[2] = 0x6A5444

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
